# feat: adding new merge table in excel
#
# Reproduces (as closely as the host object model allows):
#   - sheetView scrolled so row 13 is at the top, with D15 selected
#   - explicit custom widths on column A (~135.57 chars) and column D (~68.86 chars)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths -------------------------------------------------------
# The host quantizes ColumnWidth to steps of 1/6 of a character, so the
# inputs below are chosen as the closest reachable values to the target
# stored widths (135.5703125 and 68.85546875 character-units) rather than
# the raw target numbers themselves.
$ws.Columns.Item(1).ColumnWidth = 134.666666666667
$ws.Columns.Item(4).ColumnWidth = 68

# --- View / selection ------------------------------------------------------
# Scroll so row 13 is the top visible row, then select D15 (the author's
# last selection before saving).
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D15").Select()
